$wb = $excel.ActiveWorkbook

# --- Rename first sheet, add second sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "dces"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "downlinks"

# --- Move the "downlink" columns (D, E, F) off of sheet1 and onto the new sheet2 ---
# D (Date, duplicate) -> sheet2 column A
$ws1.Range("D1:D10").Cut($ws2.Range("A1"))
# E (Uploads, duplicate) -> sheet2 column C
$ws1.Range("E1:E10").Cut($ws2.Range("C1"))
# F (Site, duplicate) -> sheet2 column G
$ws1.Range("F1:F10").Cut($ws2.Range("G1"))

# --- Spread out the remaining columns on sheet1 ---
# B (Downloads) -> column F
$ws1.Range("B1:B10").Cut($ws1.Range("F1"))
# C (Site) -> column O
$ws1.Range("C1:C10").Cut($ws1.Range("O1"))

# --- Restore selections / active sheet to match the final layout ---
[void]$ws2.Columns.Item(4).Select()
[void]$ws1.Activate()
[void]$ws1.Range("I10").Select()
